$wb = $excel.ActiveWorkbook

# --- "About" sheet: stamp a date into C1 (serial 44316 = 4/30/2021) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 44316
$wsAbout.Range("C1").NumberFormat = "mm-dd-yy"
$wsAbout.Range("A1").Select()

# --- "RBFF" sheet: flip two fuel-fraction cells (G2: 1 -> 0, G7: 0 -> 1) ---
$wsRbff = $wb.Worksheets.Item("RBFF")
$wsRbff.Range("G2").Value = 0
$wsRbff.Range("G7").Value = 1

# The active/selected sheet moves from "About" to "RBFF"
$wsRbff.Activate()
$wsRbff.Range("M7").Select()
